$wb = $excel.ActiveWorkbook

# --- SupIm sheet: copy row 3 pattern down to rows 4-14 ---
$wsSupIm = $wb.Worksheets.Item("SupIm")
for ($i = 2; $i -le 12; $i++) {
    $row = $i + 2
    $wsSupIm.Cells.Item($row, 1).Value = $i
    $wsSupIm.Cells.Item($row, 2).Value = 0.481
    $wsSupIm.Cells.Item($row, 3).Value = 0.3
    $wsSupIm.Cells.Item($row, 4).Value = 0.207
}

# --- Demand sheet: update B3 and copy row 3 pattern down to rows 4-14 ---
$wsDemand = $wb.Worksheets.Item("Demand")
$wsDemand.Cells.Item(3, 2).Value = 224796875
for ($i = 2; $i -le 12; $i++) {
    $row = $i + 2
    $wsDemand.Cells.Item($row, 1).Value = $i
    $wsDemand.Cells.Item($row, 2).Value = 224796875
}

# --- Select Demand sheet as the active tab ---
$wsDemand.Activate()
$wsDemand.Range("E9").Select()
